# Auto-generated edit script for literature.xlsx
# Adds new evaluation tokens + updates derived statistics.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("standard-full-stats")
$ws2 = $wb.Worksheets.Item("standard-full-diffs")
$ws3 = $wb.Worksheets.Item("standard-base-stats")
$ws4 = $wb.Worksheets.Item("standard-base-diffs")
$ws5 = $wb.Worksheets.Item("nonstandard-full-stats")
$ws6 = $wb.Worksheets.Item("nonstandard-full-diffs")
$ws7 = $wb.Worksheets.Item("nonstandard-base-stats")
$ws8 = $wb.Worksheets.Item("nonstandard-base-diffs")

# ---- Update derived statistics (Precision / Recall / F1) ----
# standard-full-stats
$ws1.Range("D2").Value2 = 0.33
$ws1.Range("C3").Value2 = 0.06
$ws1.Range("D3").Value2 = 0.1
$ws1.Range("B4").Value2 = 0.84
$ws1.Range("C7").Value2 = 0.8
$ws1.Range("D7").Value2 = 0.88
$ws1.Range("C8").Value2 = 0.27
$ws1.Range("D8").Value2 = 0.33
$ws1.Range("B9").Value2 = 0.51
$ws1.Range("C9").Value2 = 0.85
$ws1.Range("D9").Value2 = 0.64

# standard-base-stats
$ws3.Range("C2").Value2 = 0.19
$ws3.Range("D2").Value2 = 0.32
$ws3.Range("C3").Value2 = 0.05
$ws3.Range("D3").Value2 = 0.07
$ws3.Range("B4").Value2 = 0.91
$ws3.Range("C4").Value2 = 0.62
$ws3.Range("D4").Value2 = 0.74
$ws3.Range("B5").Value2 = 0.64
$ws3.Range("C5").Value2 = 0.29
$ws3.Range("D5").Value2 = 0.37
$ws3.Range("B6").Value2 = 0.52
$ws3.Range("C6").Value2 = 0.87
$ws3.Range("D6").Value2 = 0.65

# nonstandard-full-stats
$ws5.Range("B4").Value2 = 0.74
$ws5.Range("D4").Value2 = 0.79
$ws5.Range("B8").Value2 = 0.58
$ws5.Range("D8").Value2 = 0.42
$ws5.Range("C9").Value2 = 0.76
$ws5.Range("D9").Value2 = 0.81

# nonstandard-base-stats
$ws7.Range("B4").Value2 = 0.82
$ws7.Range("D4").Value2 = 0.86
$ws7.Range("B5").Value2 = 0.56
$ws7.Range("D5").Value2 = 0.55
$ws7.Range("C6").Value2 = 0.79
$ws7.Range("D6").Value2 = 0.83

# ---- Append newly-tagged tokens from the updated literature sample ----
# standard-full-diffs: rows 113-133
$arr2 = New-Object 'object[,]' 21,4
$arr2[0,0] = "Jednako"
$arr2[0,1] = "Jednako"
$arr2[0,2] = "O"
$arr2[0,3] = "B-PER"
$arr2[1,0] = "!"
$arr2[1,1] = "!"
$arr2[1,2] = "O"
$arr2[1,3] = "I-PER"
$arr2[2,0] = "Brzim"
$arr2[2,1] = "Brzim"
$arr2[2,2] = "O"
$arr2[2,3] = "B-LOC"
$arr2[3,0] = "pogledom"
$arr2[3,1] = "pogledom"
$arr2[3,2] = "O"
$arr2[3,3] = "I-LOC"
$arr2[4,0] = "!"
$arr2[4,1] = "!"
$arr2[4,2] = "O"
$arr2[4,3] = "I-PER"
$arr2[5,0] = "Ja"
$arr2[5,1] = "Ja"
$arr2[5,2] = "O"
$arr2[5,3] = "B-PER"
$arr2[6,0] = "Bože"
$arr2[6,1] = "Bože"
$arr2[6,2] = "O"
$arr2[6,3] = "B-PER"
$arr2[7,0] = "!"
$arr2[7,1] = "!"
$arr2[7,2] = "O"
$arr2[7,3] = "I-PER"
$arr2[8,0] = "Za"
$arr2[8,1] = "Za"
$arr2[8,2] = "O"
$arr2[8,3] = "B-ORG"
$arr2[9,0] = "ime"
$arr2[9,1] = "ime"
$arr2[9,2] = "O"
$arr2[9,3] = "I-ORG"
$arr2[10,0] = "božje"
$arr2[10,1] = "božje"
$arr2[10,2] = "O"
$arr2[10,3] = "I-ORG"
$arr2[11,0] = "Seti"
$arr2[11,1] = "Seti"
$arr2[11,2] = "O"
$arr2[11,3] = "B-PER"
$arr2[12,0] = "Lice"
$arr2[12,1] = "Lice"
$arr2[12,2] = "O"
$arr2[12,3] = "B-PER"
$arr2[13,0] = "!"
$arr2[13,1] = "!"
$arr2[13,2] = "O"
$arr2[13,3] = "I-PER"
$arr2[14,0] = "Ja"
$arr2[14,1] = "Ja"
$arr2[14,2] = "O"
$arr2[14,3] = "B-PER"
$arr2[15,0] = "Kriv"
$arr2[15,1] = "Kriv"
$arr2[15,2] = "O"
$arr2[15,3] = "B-PER"
$arr2[16,0] = "Ti"
$arr2[16,1] = "Ti"
$arr2[16,2] = "O"
$arr2[16,3] = "B-PER"
$arr2[17,0] = "Recite"
$arr2[17,1] = "Recite"
$arr2[17,2] = "O"
$arr2[17,3] = "B-PER"
$arr2[18,0] = "Zar"
$arr2[18,1] = "Zar"
$arr2[18,2] = "O"
$arr2[18,3] = "B-PER"
$arr2[19,0] = "Recite"
$arr2[19,1] = "Recite"
$arr2[19,2] = "O"
$arr2[19,3] = "B-PER"
$arr2[20,0] = "Vaše"
$arr2[20,1] = "Vaše"
$arr2[20,2] = "O"
$arr2[20,3] = "B-ORG"
$ws2.Range("A113:D133").Value2 = $arr2

# standard-base-diffs: rows 110-130
$arr4 = New-Object 'object[,]' 21,4
$arr4[0,0] = "Jednako"
$arr4[0,1] = "Jednako"
$arr4[0,2] = "O"
$arr4[0,3] = "PER"
$arr4[1,0] = "!"
$arr4[1,1] = "!"
$arr4[1,2] = "O"
$arr4[1,3] = "PER"
$arr4[2,0] = "Brzim"
$arr4[2,1] = "Brzim"
$arr4[2,2] = "O"
$arr4[2,3] = "LOC"
$arr4[3,0] = "pogledom"
$arr4[3,1] = "pogledom"
$arr4[3,2] = "O"
$arr4[3,3] = "LOC"
$arr4[4,0] = "!"
$arr4[4,1] = "!"
$arr4[4,2] = "O"
$arr4[4,3] = "PER"
$arr4[5,0] = "Ja"
$arr4[5,1] = "Ja"
$arr4[5,2] = "O"
$arr4[5,3] = "PER"
$arr4[6,0] = "Bože"
$arr4[6,1] = "Bože"
$arr4[6,2] = "O"
$arr4[6,3] = "PER"
$arr4[7,0] = "!"
$arr4[7,1] = "!"
$arr4[7,2] = "O"
$arr4[7,3] = "PER"
$arr4[8,0] = "Za"
$arr4[8,1] = "Za"
$arr4[8,2] = "O"
$arr4[8,3] = "ORG"
$arr4[9,0] = "ime"
$arr4[9,1] = "ime"
$arr4[9,2] = "O"
$arr4[9,3] = "ORG"
$arr4[10,0] = "božje"
$arr4[10,1] = "božje"
$arr4[10,2] = "O"
$arr4[10,3] = "ORG"
$arr4[11,0] = "Seti"
$arr4[11,1] = "Seti"
$arr4[11,2] = "O"
$arr4[11,3] = "PER"
$arr4[12,0] = "Lice"
$arr4[12,1] = "Lice"
$arr4[12,2] = "O"
$arr4[12,3] = "PER"
$arr4[13,0] = "!"
$arr4[13,1] = "!"
$arr4[13,2] = "O"
$arr4[13,3] = "PER"
$arr4[14,0] = "Ja"
$arr4[14,1] = "Ja"
$arr4[14,2] = "O"
$arr4[14,3] = "PER"
$arr4[15,0] = "Kriv"
$arr4[15,1] = "Kriv"
$arr4[15,2] = "O"
$arr4[15,3] = "PER"
$arr4[16,0] = "Ti"
$arr4[16,1] = "Ti"
$arr4[16,2] = "O"
$arr4[16,3] = "PER"
$arr4[17,0] = "Recite"
$arr4[17,1] = "Recite"
$arr4[17,2] = "O"
$arr4[17,3] = "PER"
$arr4[18,0] = "Zar"
$arr4[18,1] = "Zar"
$arr4[18,2] = "O"
$arr4[18,3] = "PER"
$arr4[19,0] = "Recite"
$arr4[19,1] = "Recite"
$arr4[19,2] = "O"
$arr4[19,3] = "PER"
$arr4[20,0] = "Vaše"
$arr4[20,1] = "Vaše"
$arr4[20,2] = "O"
$arr4[20,3] = "ORG"
$ws4.Range("A110:D130").Value2 = $arr4

# nonstandard-full-diffs: rows 41-49
$arr6 = New-Object 'object[,]' 9,4
$arr6[0,0] = "Doli"
$arr6[0,1] = "Doli"
$arr6[0,2] = "B-PER"
$arr6[0,3] = "O"
$arr6[1,0] = "Doli"
$arr6[1,1] = "Doli"
$arr6[1,2] = "B-PER"
$arr6[1,3] = "O"
$arr6[2,0] = "Doli"
$arr6[2,1] = "Doli"
$arr6[2,2] = "B-PER"
$arr6[2,3] = "O"
$arr6[3,0] = "Doli"
$arr6[3,1] = "Doli"
$arr6[3,2] = "B-PER"
$arr6[3,3] = "O"
$arr6[4,0] = "Bože"
$arr6[4,1] = "Bože"
$arr6[4,2] = "O"
$arr6[4,3] = "B-PER"
$arr6[5,0] = "Doli"
$arr6[5,1] = "Doli"
$arr6[5,2] = "B-PER"
$arr6[5,3] = "O"
$arr6[6,0] = "Doli"
$arr6[6,1] = "Doli"
$arr6[6,2] = "B-PER"
$arr6[6,3] = "O"
$arr6[7,0] = "Doli"
$arr6[7,1] = "Doli"
$arr6[7,2] = "B-PER"
$arr6[7,3] = "O"
$arr6[8,0] = "Doli"
$arr6[8,1] = "Doli"
$arr6[8,2] = "B-PER"
$arr6[8,3] = "O"
$ws6.Range("A41:D49").Value2 = $arr6

# nonstandard-base-diffs: rows 37-45
$arr8 = New-Object 'object[,]' 9,4
$arr8[0,0] = "Doli"
$arr8[0,1] = "Doli"
$arr8[0,2] = "PER"
$arr8[0,3] = "O"
$arr8[1,0] = "Doli"
$arr8[1,1] = "Doli"
$arr8[1,2] = "PER"
$arr8[1,3] = "O"
$arr8[2,0] = "Doli"
$arr8[2,1] = "Doli"
$arr8[2,2] = "PER"
$arr8[2,3] = "O"
$arr8[3,0] = "Doli"
$arr8[3,1] = "Doli"
$arr8[3,2] = "PER"
$arr8[3,3] = "O"
$arr8[4,0] = "Bože"
$arr8[4,1] = "Bože"
$arr8[4,2] = "O"
$arr8[4,3] = "PER"
$arr8[5,0] = "Doli"
$arr8[5,1] = "Doli"
$arr8[5,2] = "PER"
$arr8[5,3] = "O"
$arr8[6,0] = "Doli"
$arr8[6,1] = "Doli"
$arr8[6,2] = "PER"
$arr8[6,3] = "O"
$arr8[7,0] = "Doli"
$arr8[7,1] = "Doli"
$arr8[7,2] = "PER"
$arr8[7,3] = "O"
$arr8[8,0] = "Doli"
$arr8[8,1] = "Doli"
$arr8[8,2] = "PER"
$arr8[8,3] = "O"
$ws8.Range("A37:D45").Value2 = $arr8

